# 'foaie2' - duminică 20 februarie 2022, 15:48:43 +0200 alex T460
# Update the "foaie de parcurs" (travel log) sheet from "mai 1975" to "ianuarie 2022"
# and refresh the daily travel entries accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet title / header fields -------------------------------------------------
$ws.Name = "ianuarie 2022"

$ws.Range("C5").Value = "ianuarie"
$ws.Range("C6").Value = 2022
$ws.Range("B12").Value = 89874

# Keep the print area pointed at the (renamed) sheet, same range as before.
$ws.PageSetup.PrintArea = "`$A`$1:`$F`$91"

# --- Daily travel log rows (day number in col A, Km / Loc / Observatii in B:D) --
$ws.Range("B19").Value = 152
$ws.Range("C19").Value = "Cluj-Cmp. Turzii"
$ws.Range("D19").Value = "Interes Serviciu"

$ws.Range("B21").Value = 30
$ws.Range("C21").Value = "Acasa-Birou"
$ws.Range("D21").Value = " "

$ws.Range("B23").Value = 152
$ws.Range("C23").Value = "Cluj-Cmp. Turzii"
$ws.Range("D23").Value = "Interes Serviciu"

$ws.Range("B25").Value = 156
$ws.Range("C25").Value = "Cluj-Zalau"
$ws.Range("D25").Value = "Interes Serviciu"

$ws.Range("B27").Value = 30
$ws.Range("C27").Value = "Acasa-Birou"
$ws.Range("D27").Value = " "

$ws.Range("B29:D29").Value = ""

$ws.Range("B31:D31").Value = ""

$ws.Range("B33").Value = 30
$ws.Range("C33").Value = "Acasa-Birou"
$ws.Range("D33").Value = " "

$ws.Range("B35").Value = 30
$ws.Range("C35").Value = "Acasa-Birou"
$ws.Range("D35").Value = " "

$ws.Range("B37").Value = 30
$ws.Range("C37").Value = "Acasa-Birou"
$ws.Range("D37").Value = " "

$ws.Range("B39").Value = 85
$ws.Range("C39").Value = "Cluj-Apahida"
$ws.Range("D39").Value = "Interes Serviciu"

$ws.Range("B43:D43").Value = ""

$ws.Range("B45:D45").Value = ""

$ws.Range("B47").Value = 30
$ws.Range("C47").Value = "Acasa-Birou"
$ws.Range("D47").Value = " "

$ws.Range("B49").Value = 101
$ws.Range("C49").Value = "Cluj-Dej"
$ws.Range("D49").Value = "Interes Serviciu"

$ws.Range("B51").Value = 30
$ws.Range("C51").Value = "Acasa-Birou"
$ws.Range("D51").Value = " "

$ws.Range("B53").Value = 30
$ws.Range("C53").Value = "Acasa-Birou"
$ws.Range("D53").Value = " "

$ws.Range("B57:D57").Value = ""

$ws.Range("B59:D59").Value = ""

$ws.Range("B63").Value = 30
$ws.Range("C63").Value = "Acasa-Birou"
$ws.Range("D63").Value = " "

$ws.Range("B65").Value = 30
$ws.Range("C65").Value = "Acasa-Birou"
$ws.Range("D65").Value = " "

$ws.Range("B67").Value = 356
$ws.Range("C67").Value = "Cluj-Baia-Mare"
$ws.Range("D67").Value = "Interes Serviciu"

$ws.Range("B69").Value = 257
$ws.Range("C69").Value = "Cluj-Bistrita"
$ws.Range("D69").Value = "Interes Serviciu"

$ws.Range("B71:D71").Value = ""

$ws.Range("B73:D73").Value = ""

$ws.Range("B75").Value = 356
$ws.Range("C75").Value = "Cluj-Baia-Mare"
$ws.Range("D75").Value = "Interes Serviciu"

# --- Monthly totals ---------------------------------------------------------------
$ws.Range("B76").Value = 2046
$ws.Range("B77").Value = 91920
